$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number, date range)
$ws.Range("A8").Value = "Volume 31   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# Row 14
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -50
$ws.Range("L14").Value = -25
$ws.Range("N14").Value = -83.333333333333

# Row 15
$ws.Range("G14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("G14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -31.578947368421
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = -7.142857142857
$ws.Range("N15").Value = -76.785714285714

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -21.739130434782
$ws.Range("I16").Value = 193
$ws.Range("J16").Value = 237
$ws.Range("K16").Value = -18.565400843881
$ws.Range("L16").Value = -13.452914798206
$ws.Range("M16").Value = -22.8
$ws.Range("N16").Value = -85.108024691358

# Row 17
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").Value = 73
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 43.137254901960
$ws.Range("I17").Value = 467
$ws.Range("J17").Value = 423
$ws.Range("K17").Value = 10.401891252955
$ws.Range("L17").Value = 27.945205479452
$ws.Range("M17").Value = 80.308880308880
$ws.Range("N17").Value = -33.190271816881

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -39.130434782608
$ws.Range("I18").Value = 96
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = -28.358208955223
$ws.Range("L18").Value = -13.513513513513
$ws.Range("M18").Value = -17.241379310344
$ws.Range("N18").Value = -93.048515568428

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 4.761904761904
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 365
$ws.Range("K19").Value = -10.958904109589
$ws.Range("L19").Value = -5.523255813953
$ws.Range("M19").Value = 77.595628415300
$ws.Range("N19").Value = -43.280977312390

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 9.523809523809
$ws.Range("I20").Value = 171
$ws.Range("J20").Value = 139
$ws.Range("K20").Value = 23.021582733812
$ws.Range("L20").Value = 106.024096385542
$ws.Range("M20").Value = 122.077922077922
$ws.Range("N20").Value = -69.680851063829

# Row 21
$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 62.5
$ws.Range("F21").Value = 176
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = 4.761904761904
$ws.Range("I21").Value = 1271
$ws.Range("J21").Value = 1325
$ws.Range("K21").Value = -4.075471698113
$ws.Range("L21").Value = 10.907504363001
$ws.Range("M21").Value = 40.132304299889
$ws.Range("N21").Value = -72.399565689468

# Row 22
$ws.Range("F22").Value = 1
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -29.411764705882

# Row 23
$ws.Range("D14").Copy($ws.Range("C23"))
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("L23").Value = -33.333333333333

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -53.846153846153
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -29.838709677419
$ws.Range("I24").Value = 626
$ws.Range("J24").Value = 851
$ws.Range("K24").Value = -26.439482961222
$ws.Range("L24").Value = 5.743243243243
$ws.Range("M24").Value = 6.101694915254

# Row 25
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 30
$ws.Range("F25").Value = 97
$ws.Range("G25").Value = 84
$ws.Range("H25").Value = 15.476190476190
$ws.Range("I25").Value = 765
$ws.Range("J25").Value = 638
$ws.Range("K25").Value = 19.905956112852
$ws.Range("L25").Value = 47.969052224371
$ws.Range("M25").Value = -1.290322580645

# Row 26
$ws.Range("G14").Copy($ws.Range("C26"))
$ws.Range("C26").Value = 2
$ws.Range("G14").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("K14").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = -28.571428571428
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 32
$ws.Range("K26").Value = 9.375
$ws.Range("L26").Value = 52.173913043478

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = -7.692307692307

# Row 28
$ws.Range("D14").Copy($ws.Range("C28"))
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 66.666666666666
$ws.Range("L28").Value = 25.925925925925
$ws.Range("M28").Value = -2.857142857142
$ws.Range("N28").Value = -68.224299065420

# Row 29
$ws.Range("D14").Copy($ws.Range("C29"))
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -33.333333333333
$ws.Range("L29").Value = -7.692307692307
$ws.Range("M29").Value = -17.241379310344
$ws.Range("N29").Value = -76.237623762376
